$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: new talk entry ("HSC Year 3 Weak Lensing Cosmology Results")
$ws.Range("A27").Value = "HSC Year 3 Weak Lensing Cosmology Results"
$ws.Range("C27").Value = 2023
$ws.Range("D27").Value = 4
$ws.Range("F27").Value = "HSC webinar"
$ws.Range("J27").Value = 1
$ws.Range("K27").Value = 1

$oCell = $ws.Range("O27")
$ws.Hyperlinks.Add($oCell, "https://hsc-release.mtk.nao.ac.jp/doc/index.php/wly3/")
$oCell.Font.Underline = $true
$oCell.Font.ThemeColor = 11

$ws.Range("Q27").Value = "We presented our HSC-Y3 weak lensing cosmology results on webinar"
$ws.Range("R27").Value = 1

# Move the active selection like the saved workbook did
$ws.Range("B28").Select()
